$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook identity: sheet (tab) name ---
$ws.Name = "B08L7BD4C6"

# --- Column A, rows 1-100: replacement keyword list ---
$values = New-Object "object[,]" 100,1
$values[0,0] = "black sport"
$values[1,0] = "clothing set"
$values[2,0] = "seamless legging"
$values[3,0] = "seamless crop top"
$values[4,0] = "legging suit"
$values[5,0] = "legging top set"
$values[6,0] = "green apparel"
$values[7,0] = "small crop top"
$values[8,0] = "black compression"
$values[9,0] = "winter activewear"
$values[10,0] = "winter blue"
$values[11,0] = "winter tracksuit"
$values[12,0] = "pink work shirt"
$values[13,0] = "work out outfit"
$values[14,0] = "winter outfit"
$values[15,0] = "2 piece matching set"
$values[16,0] = "2 piece long sleeve"
$values[17,0] = "2 piece workout set"
$values[18,0] = "black pant suit"
$values[19,0] = "black pant set"
$values[20,0] = "black workout shirt"
$values[21,0] = "black workout top"
$values[22,0] = "seamless top"
$values[23,0] = "pink tight shirt"
$values[24,0] = "pink crop shirt"
$values[25,0] = "green high top"
$values[26,0] = "work out crop"
$values[27,0] = "two piece legging set"
$values[28,0] = "seamless workout top"
$values[29,0] = "seamless workout"
$values[30,0] = "black sport shirt"
$values[31,0] = "winter workout"
$values[32,0] = "long black top"
$values[33,0] = "tight long sleeve workout shirt"
$values[34,0] = "workout tracksuit"
$values[35,0] = "workout pant"
$values[36,0] = "2 piece jogger"
$values[37,0] = "2 piece woman"
$values[38,0] = "athletic legging"
$values[39,0] = "black tight top"
$values[40,0] = "black jogger"
$values[41,0] = "black long crop top"
$values[42,0] = "crop top woman"
$values[43,0] = "green black top"
$values[44,0] = "long sleeve winter shirt"
$values[45,0] = "2 piece pant set"
$values[46,0] = "2 piece pant suit"
$values[47,0] = "long sleeve work out shirt"
$values[48,0] = "long sleeve set"
$values[49,0] = "black top woman"
$values[50,0] = "long sleeve outfit"
$values[51,0] = "crop top pant set"
$values[52,0] = "long sleeve athletic top"
$values[53,0] = "long sleeve two piece set"
$values[54,0] = "pink workout set"
$values[55,0] = "black jogger set"
$values[56,0] = "tight workout top"
$values[57,0] = "long sleeve clothing"
$values[58,0] = "black athletic crop top"
$values[59,0] = "small black crop top"
$values[60,0] = "black athletic sleeve"
$values[61,0] = "matching legging set"
$values[62,0] = "long sleeve workout top"
$values[63,0] = "long sleeve crop top workout"
$values[64,0] = "pink long sleeve crop"
$values[65,0] = "black long sleeve crop"
$values[66,0] = "green long sleeve crop"
$values[67,0] = "long sleeve crop top set"
$values[68,0] = "long sleeve work"
$values[69,0] = "top pant set"
$values[70,0] = "pink outfit set"
$values[71,0] = "black crop"
$values[72,0] = "long black"
$values[73,0] = "work outfit"
$values[74,0] = "piece work"
$values[75,0] = "set two"
$values[76,0] = "fitness pant"
$values[77,0] = "winter shirt"
$values[78,0] = "workout tight"
$values[79,0] = "high top woman"
$values[80,0] = "matching pant set"
$values[81,0] = "gym tight"
$values[82,0] = "legging outfit"
$values[83,0] = "long green"
$values[84,0] = "crop legging"
$values[85,0] = "work set"
$values[86,0] = "seamless long sleeve workout top"
$values[87,0] = "two apparel"
$values[88,0] = "tracksuit pant"
$values[89,0] = "work out pink"
$values[90,0] = "workout track"
$values[91,0] = "black high top"
$values[92,0] = "fashionable top"
$values[93,0] = "winter top"
$values[94,0] = "2 small"
$values[95,0] = "seamless crop"
$values[96,0] = "yoga with"
$values[97,0] = "workout jogger"
$values[98,0] = "athletic wear"
$values[99,0] = "sport legging"
$ws.Range("A1:A100").Value = $values

# --- Selection moves to E13 ---
$ws.Range("E13").Select() | Out-Null

